$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 96, pushing the existing 96:128
# data block down to 98:130 (matches dimension growing from T128 to T130).
$ws.Range("A96:A97").EntireRow.Insert()

# --- New row 96 ---
$ws.Range("A96").Value = 10
$ws.Range("B96").Value = "Vega Modelo de Temuco"
$ws.Range("C96").Value = "La Araucanía"
$ws.Range("D96").Value = Get-Date -Year 2022 -Month 12 -Day 29 -Hour 0 -Minute 0 -Second 0
$ws.Range("E96").Value = 9
$ws.Range("F96").Value = "Fruta"
$ws.Range("G96").Value = 100101
$ws.Range("H96").Value = "Berries"
$ws.Range("I96").Value = 100101001
$ws.Range("J96").Value = "Arándano (blue)"
$ws.Range("K96").Value = "Sin especificar"
$ws.Range("L96").Value = "Primera"
$ws.Range("M96").Value = 230
$ws.Range("N96").Value = 1800
$ws.Range("O96").Value = 2000
$ws.Range("P96").Value = 1870
$ws.Range("Q96").Value = "$/kilo"
$ws.Range("R96").Value = "Provincia de Cardenal Caro"
$ws.Range("S96").Value = 1870
$ws.Range("T96").Value = 1

# --- New row 97 ---
$ws.Range("A97").Value = 10
$ws.Range("B97").Value = "Vega Modelo de Temuco"
$ws.Range("C97").Value = "La Araucanía"
$ws.Range("D97").Value = Get-Date -Year 2022 -Month 12 -Day 29 -Hour 0 -Minute 0 -Second 0
$ws.Range("E97").Value = 9
$ws.Range("F97").Value = "Fruta"
$ws.Range("G97").Value = 100101
$ws.Range("H97").Value = "Berries"
$ws.Range("I97").Value = 100101001
$ws.Range("J97").Value = "Arándano (blue)"
$ws.Range("K97").Value = "Sin especificar"
$ws.Range("L97").Value = "Segunda"
$ws.Range("M97").Value = 110
$ws.Range("N97").Value = 1500
$ws.Range("O97").Value = 1500
$ws.Range("P97").Value = 1500
$ws.Range("Q97").Value = "$/kilo"
$ws.Range("R97").Value = "Provincia de Cardenal Caro"
$ws.Range("S97").Value = 1500
$ws.Range("T97").Value = 1
